# Revert "Removed knock outs in iFerment"
# Restores original flux values (column B) in Sheet1 that were zeroed out
# by the prior commit, reproducing the simulation results exactly.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B13").Value = [double]"0.01647178444331563"
$ws.Range("B14").Value = [double]"-0.03924758515505806"
$ws.Range("B15").Value = [double]"-0.03924758515505805"
$ws.Range("B16").Value = [double]"-0.003457041179461431"
$ws.Range("B17").Value = [double]"-0.003457041179461432"
$ws.Range("B18").Value = [double]"-0.05226232841891227"
$ws.Range("B19").Value = [double]"-0.9999999999999983"
$ws.Range("B23").Value = [double]"0.9682765632943552"
$ws.Range("B25").Value = [double]"0.9125571936959815"
$ws.Range("B26").Value = [double]"0.9125571936959815"
$ws.Range("B27").Value = [double]"0.9125571936959814"
$ws.Range("B28").Value = [double]"1.596542958820536"
$ws.Range("B29").Value = [double]"-1.596542958820536"
$ws.Range("B30").Value = [double]"-1.436705643111334"
$ws.Range("B31").Value = [double]"1.436705643111334"
$ws.Range("B32").Value = [double]"0"
$ws.Range("B40").Value = [double]"0.03538383324860186"
$ws.Range("B41").Value = [double]"0.05246568378240964"
$ws.Range("B42").Value = [double]"0"
$ws.Range("B43").Value = [double]"-0.06893746822572527"
$ws.Range("B45").Value = [double]"0.4471784443314684"
$ws.Range("B47").Value = [double]"0"
$ws.Range("B48").Value = [double]"0"
$ws.Range("B49").Value = [double]"0"
$ws.Range("B50").Value = [double]"0"
$ws.Range("B55").Value = [double]"0"
$ws.Range("B56").Value = [double]"0"
$ws.Range("B57").Value = [double]"0"
$ws.Range("B58").Value = [double]"0"
$ws.Range("B63").Value = [double]"0"
$ws.Range("B64").Value = [double]"0"
$ws.Range("B65").Value = [double]"0"
$ws.Range("B66").Value = [double]"0"
$ws.Range("B69").Value = [double]"0"
$ws.Range("B70").Value = [double]"0"
$ws.Range("B104").Value = [double]"9.088154550076236"
$ws.Range("B106").Value = [double]"0"
$ws.Range("B108").Value = [double]"0"
$ws.Range("B109").Value = [double]"0"
$ws.Range("B118").Value = [double]"0.4471784443314684"
$ws.Range("B119").Value = [double]"3.607320793085909"
$ws.Range("B120").Value = [double]"-3.370615149974573"
$ws.Range("B121").Value = [double]"3.370615149974573"
$ws.Range("B122").Value = [double]"0.9285205897305565"
$ws.Range("B123").Value = [double]"-0.6729028978139283"
$ws.Range("B124").Value = [double]"0.6729028978139283"
$ws.Range("B127").Value = [double]"0.1427554651753939"
$ws.Range("B129").Value = [double]"0.1427554651753939"
$ws.Range("B130").Value = [double]"0.04514489069649202"
$ws.Range("B133").Value = [double]"0"
$ws.Range("B136").Value = [double]"-7.803729132524266e-16"
$ws.Range("B137").Value = [double]"-6.767463141840349"
$ws.Range("B138").Value = [double]"-1.051753940010165"
$ws.Range("B139").Value = [double]"1.051753940010165"
$ws.Range("B140").Value = [double]"-0.02562277580071169"
$ws.Range("B141").Value = [double]"0.02562277580071168"
$ws.Range("B142").Value = [double]"0.1708185053380779"
$ws.Range("B143").Value = [double]"0.5783426537874925"
$ws.Range("B144").Value = [double]"-0.1049313675648182"
$ws.Range("B145").Value = [double]"0.1049313675648182"
$ws.Range("B147").Value = [double]"0.1220132180986271"
$ws.Range("B148").Value = [double]"0.1220132180986271"
$ws.Range("B151").Value = [double]"0.03538383324860186"
$ws.Range("B153").Value = [double]"0"
$ws.Range("B165").Value = [double]"0"
$ws.Range("B175").Value = [double]"0.8931367564819547"
$ws.Range("B188").Value = [double]"0.7851550584646669"
$ws.Range("B192").Value = [double]"-4.392475851550575"
$ws.Range("B199").Value = [double]"-0.3244534824605992"
$ws.Range("B200").Value = [double]"-0.0696492119979665"
$ws.Range("B210").Value = [double]"0.9999999999999982"
$ws.Range("B241").Value = [double]"0.04514489069649202"
$ws.Range("B242").Value = [double]"9.088154550076235"
